$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.580.17'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.960.93'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'243.58"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  +2.17%  '
$ws.Range('D7').Value = "'60.35"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.27%  '
$ws.Range('D8').Value = "'1.00"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'0.378"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.43%  '
$ws.Range('D10').Value = "'0.0787"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('D12').Value = "'14.16"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.83%  '
$ws.Range('D13').Value = "'0.844"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('D14').Value = '2.245.05'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').Value = "'21.62"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = "'5.28"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = '1.954.31'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '36.501.33'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').Value = "'69.26"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').Value = "'229.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').Value = "'5.08"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = "'2.44"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('D25').Value = "'2.36"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.17%  '
$ws.Range('D26').Value = "'0.145"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.94%  '
$ws.Range('D27').Value = "'9.14"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.17%  '
$ws.Range('D28').Value = "'161.13"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').Value = "'19.30"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  +19.18%  '
$ws.Range('D31').Value = "'0.121"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').Value = "'4.79"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').Value = "'0.0613"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').Value = "'4.48"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.56%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = "'2.27"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').Value = "'1.77"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('D39').Value = "'5.41"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.19%  '
$ws.Range('D40').Value = "'0.0968"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').Value = "'15.85"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.97%  '
$ws.Range('D45').Value = '1.363.45'
$ws.Range('E45').Value = '  +1.37%  '
$ws.Range('D46').Value = "'88.50"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.65%  '
$ws.Range('D47').Value = "'1.03"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('D48').Value = "'7.22"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').Value = "'45.91"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.89%  '
$ws.Range('D51').Value = '2.140.93'
$ws.Range('E51').Value = '  +0.64%  '
